# Updating FoodLog file 06/05/2018
# Appends 7 new days of food-log data (30 Apr 2018 - 06 May 2018) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Prepare rows 120-126 with the same cell formatting as row 119
#    (this keeps Column A/B/C..H tied to the existing style records
#    instead of Excel's "General" default).
# ------------------------------------------------------------------
$ws.Range("A119:J119").Copy()
$ws.Range("A120:J126").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column A for the new rows gets its own date number format (separate
# from the existing custom "d/mm/yyyy;@" style used by the rest of the
# date column).
$ws.Range("A120:A126").NumberFormat = "m/d/yyyy"

# ------------------------------------------------------------------
# 2. New data values
# ------------------------------------------------------------------
$rows = @(
    @(120, 43220, 1884, 57, 27, 162, 3415, 178, 3000),
    @(121, 43221, 1679, 54, 25, 151, 2744, 144, 3000),
    @(122, 43222, 1706, 66, 16, 176, 2429, 100, 3000),
    @(123, 43223, 1227, 38, 19, 132, 2382,  90, 2500),
    @(124, 43224, 1370, 43, 15, 167, 2791,  74, 3500),
    @(125, 43225, 2274, 63, 41, 257, 3950, 182, 3250),
    @(126, 43226, 1637, 49, 24, 174, 1966, 121, 3000)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]   # A - Date
    $ws.Cells.Item($row, 2).Value = $r[2]   # B - Calories In
    $ws.Cells.Item($row, 3).Value = $r[3]   # C
    $ws.Cells.Item($row, 4).Value = $r[4]   # D
    $ws.Cells.Item($row, 5).Value = $r[5]   # E
    $ws.Cells.Item($row, 6).Value = $r[6]   # F
    $ws.Cells.Item($row, 7).Value = $r[7]   # G
    $ws.Cells.Item($row, 8).Value = $r[8]   # H
}

# ------------------------------------------------------------------
# 3. Extend the I/J formulas down through row 126
# ------------------------------------------------------------------
$ws.Range("I120:I126").Formula = "=IF(H120>=2200,""Yes"",""No"")"
$ws.Range("J120:J126").Formula = "=IF(B120<=1800,""Yes"",""No"")"

# ------------------------------------------------------------------
# 4. Update the view: scroll/selection now rests on G124, and the old
#    frozen "topLeftCell" (A81) should no longer be pinned.
# ------------------------------------------------------------------
$ws.Range("G124").Select()
